$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.244.59'
$ws.Range('E2').Value = '  +0.43%  '

$ws.Range('D3').Value = '3.856.19'
$ws.Range('E3').Value = '  +0.03%  '

$ws.Range('E4').Value = '  +0.18%  '

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '464.53'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +9.49%  '

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '147.78'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +13.53%  '

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.631'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  +3.41%  '

$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.998'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.01%  '

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.750'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +3.46%  '

$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.155'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -2.79%  '

$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '0.0000311'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -8.32%  '

$ws.Range('E12').Value = '  +7.90%  '

$ws.Range('E13').Value = '  +1.43%  '

$ws.Range('D14').Value = '4.483.80'
$ws.Range('E14').Value = '  +0.29%  '

$ws.Range('E15').Value = '  -7.09%  '

$ws.Range('D16').Value = '3.859.05'
$ws.Range('E16').Value = '  +0.34%  '

$ws.Range('E17').Value = '  -0.11%  '

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '20.02'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +0.36%  '

$ws.Range('E19').Value = '  +7.48%  '

$ws.Range('D20').Value = '67.459.12'
$ws.Range('E20').Value = '  +0.37%  '

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '431.49'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +4.13%  '

$ws.Range('E22').Value = '  -0.86%  '

$ws.Range('E23').Value = '  +7.69%  '

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '88.50'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +4.96%  '

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '3.55'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +9.23%  '

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '10.49'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +13.91%  '

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '37.58'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  -0.14%  '

$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '10.20'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +2.00%  '

$ws.Range('E29').Value = '  +4.03%  '

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '745.07'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +1.47%  '

$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.134'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +9.59%  '

$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '13.74'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +4.14%  '

$ws.Range('E33').Value = '  +2.09%  '

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '43.29'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +11.78%  '

$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '0.162'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +6.66%  '

$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '57.22'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +3.20%  '

$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.999'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.07%  '

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '5.55'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +0.84%  '

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.0479'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +3.26%  '

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.352'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +11.09%  '

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.93'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.14%  '

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.63'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +14.30%  '

$ws.Range('E43').Value = '  +5.09%  '

$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +0.26%  '

$ws.Range('B45').Value = 'PEPE'
$ws.Range('C45').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D45').Value = '0.0₃0672'
$ws.Range('E45').Value = '  -9.20%  '

$ws.Range('B46').Value = 'LidoDAOToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '3.44'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +2.31%  '

$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '3.27'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +4.55%  '

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '2.76'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +7.65%  '

$ws.Range('E49').Value = '  +3.48%  '

$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.90'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +2.87%  '

$ws.Range('B51').Value = 'Monero'
$ws.Range('C51').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '144.56'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +3.47%  '
